# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Update header labels on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" (last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse existing header style (bold, centered, bordered) from "Weekly Quantity"!A1:B1
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Reuse existing date style from "Weekly Quantity"!A2 for the "ds" column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$wsForecast.Range("A2").Value = 44934.99999999999
$wsForecast.Range("B2").Value = 82
$wsForecast.Range("C2").Value = 32.8106346068572
$wsForecast.Range("D2").Value = 130.3917844971505

$wsForecast.Range("A3").Value = 44941.99999999999
$wsForecast.Range("B3").Value = 65
$wsForecast.Range("C3").Value = 15.29005933314723
$wsForecast.Range("D3").Value = 111.4934530751888

$wsForecast.Range("A4").Value = 44955.99999999999
$wsForecast.Range("B4").Value = 29
$wsForecast.Range("C4").Value = -19.1631649684508
$wsForecast.Range("D4").Value = 76.45178743741521

$wsForecast.Range("A5").Value = 44962.99999999999
$wsForecast.Range("B5").Value = 12
$wsForecast.Range("C5").Value = -34.55474300199985
$wsForecast.Range("D5").Value = 60.50352816750845

$wsForecast.Range("A6").Value = 44969.99999999999
$wsForecast.Range("B6").Value = 0
$wsForecast.Range("C6").Value = -56.02516963328503
$wsForecast.Range("D6").Value = 43.3965534444151

$wsForecast.Range("A7").Value = 44976.99999999999
$wsForecast.Range("B7").Value = 0
$wsForecast.Range("C7").Value = -73.12172271194886
$wsForecast.Range("D7").Value = 24.6483212102739

$wsForecast.Range("A8").Value = 44983.99999999999
$wsForecast.Range("B8").Value = 0
$wsForecast.Range("C8").Value = -89.65413981804406
$wsForecast.Range("D8").Value = 11.20816711963501

$wsForecast.Range("A9").Value = 44990.99999999999
$wsForecast.Range("B9").Value = 0
$wsForecast.Range("C9").Value = -105.533067378588
$wsForecast.Range("D9").Value = -9.316255081353717

$wsForecast.Range("A10").Value = 44997.99999999999
$wsForecast.Range("B10").Value = 0
$wsForecast.Range("C10").Value = -125.1109432205859
$wsForecast.Range("D10").Value = -28.45407595662616

$wsForecast.Range("A11").Value = 45004.99999999999
$wsForecast.Range("B11").Value = 0
$wsForecast.Range("C11").Value = -144.8031826139323
$wsForecast.Range("D11").Value = -43.67925063799549

$wsForecast.Range("A12").Value = 45011.99999999999
$wsForecast.Range("B12").Value = 0
$wsForecast.Range("C12").Value = -162.4401600478779
$wsForecast.Range("D12").Value = -59.70512571118288

# Select A1 on the new sheet, matching target selection state
$wsForecast.Range("A1").Select()
